# Update scripts with new TPM values (Fam3c-Lifr sheet)
#
# A re-quantification of the TPM source data changed the ligand average
# expression for the "ECs" sending cluster and the receptor average
# expression for the "ECs" target cluster. Every other figure in this table
# (totals, derived-specificity scores and edge weights) is a deterministic
# function of those two numbers, so the dependent columns are recomputed
# below for each of the 9 rows (sending cluster x target cluster) and only
# the cells whose value actually changes are written back to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated per-cluster ligand average expression (column G), keyed by the
# "Sending cluster" value in column A. Only ECs actually changed with the
# new TPM quantification; the others are listed so the specificity totals
# below are computed against the full, consistent set of values.
$ligandAvg = @{
    "ECs"   = 6.480005333333334
    "FAPs"  = 14.13955433333333
    "MuSCs" = 4.063282999999999
}

# Updated per-cluster receptor average expression (column M), keyed by the
# "Target cluster" value in column D. Only ECs actually changed.
$receptorAvg = @{
    "ECs"   = 34.53319033333333
    "FAPs"  = 80.77474466666666
    "MuSCs" = 24.96420533333334
}

$ligandCells = @{ "ECs" = 3; "FAPs" = 3; "MuSCs" = 3 }
$receptorCells = @{ "ECs" = 3; "FAPs" = 3; "MuSCs" = 3 }

$ligandTotal = ($ligandAvg.Values | Measure-Object -Sum).Sum
$receptorTotal = ($receptorAvg.Values | Measure-Object -Sum).Sum

# Only the clusters below actually had their raw average expression value
# changed by the new TPM numbers; rows that reference only unaffected
# clusters keep columns G/H/M/N/Q/R exactly as they were. Every row's
# derived-specificity columns (I/J, O/P) and the edge-specificity columns
# (S/T) still move because they are normalised against the column totals,
# which shifted once any cluster's average changed.
$changedLigandClusters = @("ECs")
$changedReceptorClusters = @("ECs")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $sending = $ws.Range("A$r").Value()
    $target = $ws.Range("D$r").Value()

    $ligandChanged = $changedLigandClusters -contains $sending
    $receptorChanged = $changedReceptorClusters -contains $target

    $g = $ligandAvg[$sending]
    $h = $g * $ligandCells[$sending]
    $i = $g / $ligandTotal

    $m = $receptorAvg[$target]
    $n = $m * $receptorCells[$target]
    $o = $m / $receptorTotal

    if ($ligandChanged) {
        $ws.Range("G$r").Value = $g
        $ws.Range("H$r").Value = $h
    }

    # Ligand derived specificity depends on every sending-cluster average,
    # so it is recomputed for every row whenever any of them changes.
    $ws.Range("I$r").Value = $i
    $ws.Range("J$r").Value = $i

    if ($receptorChanged) {
        $ws.Range("M$r").Value = $m
        $ws.Range("N$r").Value = $n
    }

    # Receptor derived specificity likewise depends on every target-cluster
    # average.
    $ws.Range("O$r").Value = $o
    $ws.Range("P$r").Value = $o

    # Edge average/total expression weight (Q/R) are simple products of the
    # raw (non-normalised) averages/totals, so they only move when the raw
    # ligand or receptor value on this row actually changed.
    if ($ligandChanged -or $receptorChanged) {
        $q = $g * $m
        $rr = $h * $n
        $ws.Range("Q$r").Value = $q
        $ws.Range("R$r").Value = $rr
    }

    # Edge derived specificity (S/T) is the product of the two normalised
    # specificity scores, so it moves whenever either normalised score
    # moves - i.e. on every row, since the totals shifted.
    $s = $i * $o
    $ws.Range("S$r").Value = $s
    $ws.Range("T$r").Value = $s
}
